# Weekly update: a new weekly price record is inserted as row 520
# (Feria Lagunitas de Puerto Montt - Coliflor), pushing the existing
# historical rows 520-537 down to 521-538.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 520; this shifts rows 520:537 down to 521:538
# and keeps the dimension/used-range in sync (A1:R537 -> A1:R538).
$ws.Rows.Item(520).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(520, 1).Value = 4
$ws.Cells.Item(520, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(520, 3).Value = "Los Lagos"
$ws.Cells.Item(520, 4).Value = 45075
$ws.Cells.Item(520, 5).Value = 10
$ws.Cells.Item(520, 6).Value = 100112008
$ws.Cells.Item(520, 7).Value = "Coliflor"
$ws.Cells.Item(520, 8).Value = "Sin especificar"
$ws.Cells.Item(520, 9).Value = "Primera"
$ws.Cells.Item(520, 10).Value = 250
$ws.Cells.Item(520, 11).Value = 1700
$ws.Cells.Item(520, 12).Value = 1700
$ws.Cells.Item(520, 13).Value = 1700
$ws.Cells.Item(520, 14).Value = "$/unidad"
$ws.Cells.Item(520, 15).Value = "Región Metropolitana"
$ws.Cells.Item(520, 16).Value = 1700
$ws.Cells.Item(520, 17).Value = 1
$ws.Cells.Item(520, 18).Value = "Hortaliza"
